# "Updated symbol list ... with GitHub Actions" — refresh the scraped coin
# price/volume snapshot in column D (Price) and a couple of Volume(1h)
# strings in column E.
#
# The sheet stores Price as TEXT (not Number) even though the values look
# numeric — that's how the upstream scraper wrote the workbook. A plain
# `Range.Value = "246.77"` would get auto-coerced to a Number by Excel, so
# instead we stage each value in a scratch cell that's explicitly formatted
# as Text ("@"), then Copy/PasteSpecial *values only* into the destination.
# That keeps the destination cell's own style/format untouched (no stray
# NumberFormat or quote-prefix residue) while still landing the value as a
# genuine Text cell, matching the original file.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$helper = $ws.Range("ZZ1")

function Set-TextPrice([string]$cellAddr, [string]$newValue) {
    $helper.NumberFormat = "@"
    $helper.Value = $newValue
    $helper.Copy()
    $ws.Range($cellAddr).PasteSpecial(-4163)  # xlPasteValues
}

$priceUpdates = @{
    "D2"  = "246.77"
    "D3"  = "22.68"
    "D4"  = "5.253"
    "D5"  = "0.05696"
    "D7"  = "6.291"
    "D8"  = "0.8097"
    "D9"  = "0.8689"
    "D10" = "0.1426"
    "D11" = "0.07354"
    "D12" = "0.03028"
    "D14" = "0.09391"
    "D15" = "3.869"
    "D16" = "0.001575"
    "D17" = "0.04786"
    "D18" = "0.0005851"
    "D19" = "0.006081"
    "D20" = "0.005025"
    "D21" = "0.0009969"
    "D40" = "0.03928"
    "D41" = "0.006789"
    "D42" = "0.1068"
    "D43" = "0.002679"
    "D44" = "0.007504"
    "D48" = "0.1928"
}

foreach ($cellAddr in $priceUpdates.Keys) {
    Set-TextPrice $cellAddr $priceUpdates[$cellAddr]
}

# Volume(1h) label text tweaks (plain, non-numeric strings — no coercion risk).
$ws.Range("E18").Value = "17OneONE"
$ws.Range("E48").Value = "47BOLOBOLOWorstin24h"

$helper.Clear()
$ws.Range("A1").Select() | Out-Null
